# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
# For every cell listed below we overwrite the old scraped value with the new one.
# Columns: B=Coin name, C=Link, D=Price, E=Volume(1h) change.
#
# Numeric-looking price strings (column D) are forced to Text format first so Excel
# does not silently reinterpret e.g. "4.30" as the number 4.3, which would not match
# the original text content of the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "36.162.88"
$ws.Range("E2").Value = "  -3.68%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.967.78"
$ws.Range("E3").Value = "  -2.73%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.48"
$ws.Range("E5").Value = "  -12.82%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.595"
$ws.Range("E6").Value = "  -3.92%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.57"
$ws.Range("E8").Value = "  -4.62%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  -5.54%  "

# Row 10: OKB
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.11"
$ws.Range("E10").Value = "  +0.28%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("E11").Value = "  -4.77%  "

# Row 12: TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0978"
$ws.Range("E12").Value = "  -3.99%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.255.43"
$ws.Range("E13").Value = "  -3.01%  "

# Row 14: Chainlink
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.71"
$ws.Range("E14").Value = "  -5.06%  "

# Row 15: Avalanche
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.67"
$ws.Range("E15").Value = "  -5.95%  "

# Row 16: Polygon
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.749"
$ws.Range("E16").Value = "  -7.40%  "

# Row 17: Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.96"
$ws.Range("E17").Value = "  -5.80%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "1.969.42"
$ws.Range("E18").Value = "  -3.27%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "36.103.09"
$ws.Range("E19").Value = "  -3.39%  "

# Row 20: Litecoin
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.16"
$ws.Range("E20").Value = "  -3.97%  "

# Row 21: ShibaInu
$ws.Range("D21").Value = "0.0₃0801"
$ws.Range("E21").Value = "  -5.36%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.01"
$ws.Range("E22").Value = "  -3.34%  "

# Row 23: BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.77"
$ws.Range("E23").Value = "  -3.93%  "

# Row 24: Dai
$ws.Range("E24").Value = "  +0.16%  "

# Row 25: Toncoin
$ws.Range("E25").Value = "  +0.34%  "

# Row 26: PancakeSwap
$ws.Range("E26").Value = "  -14.76%  "

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.34"
$ws.Range("E27").Value = "  -3.02%  "

# Row 28: Cosmos
$ws.Range("E28").Value = "  -6.21%  "

# Row 29: EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.56"
$ws.Range("E29").Value = "  -5.94%  "

# Row 30: ImmutableX
$ws.Range("E30").Value = "  -2.74%  "

# Row 31: Kaspa
$ws.Range("E31").Value = "  -6.24%  "

# Row 32: Stellar
$ws.Range("E32").Value = "  -3.95%  "

# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.30"
$ws.Range("E33").Value = "  -7.39%  "

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0597"
$ws.Range("E34").Value = "  -8.74%  "

# Row 35: InternetComputer(DFINITY)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.20"
$ws.Range("E35").Value = "  -7.28%  "

# Row 36: LidoDAOToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.29"
$ws.Range("E36").Value = "  -4.07%  "

# Row 37: BinanceUSD
$ws.Range("E37").Value = "  -0.23%  "

# Row 38: WEMIXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -2.16%  "

# Row 39: RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.20"
$ws.Range("E39").Value = "  -4.17%  "

# Row 40: THORChain
$ws.Range("E40").Value = "  -2.15%  "

# Row 41: HuobiToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  -1.06%  "

# Row 42: Maker
$ws.Range("D42").Value = "1.418.81"
$ws.Range("E42").Value = "  +0.83%  "

# Row 43: Cronos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0877"
$ws.Range("E43").Value = "  -7.56%  "

# Row 44: VeChain
$ws.Range("E44").Value = "  -7.44%  "

# Row 45: TrustWalletToken
$ws.Range("E45").Value = "  -12.69%  "

# Row 46: Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.88"
$ws.Range("E46").Value = "  -4.58%  "

# Row 47: ARBITRUM
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.978"
$ws.Range("E47").Value = "  -5.50%  "

# Row 48: MXToken
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "14.56"
$ws.Range("E48").Value = "  -7.28%  "

# Row 49: InjectiveProtocol
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.85"
$ws.Range("E49").Value = "  -0.85%  "

# Row 50: FraxShare
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.68"
$ws.Range("E50").Value = "  -6.03%  "

# Row 51: FTXToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.56"
$ws.Range("E51").Value = "  +12.65%  "
